$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Which formula represents the Pythagorean theorem? -> Which HTTP status code means "Not Found"?
$ws.Range("A2").Value = 'Which HTTP status code means "Not Found"?'
$ws.Range("B2").Value = "Select the correct numeric status code"
$ws.Range("C2").Value = "Web Dev"
$ws.Range("D2").Value = "single_choice"
$ws.Range("E2").Value = '["200","301","404"]'
$ws.Range("F2").Value = '["404"]'
$ws.Range("H2").Value = "http,status"

# Row 3: Explain the Pythagorean theorem -> Explain React useState/useEffect
$ws.Range("A3").Value = "Explain the difference between React useState and useEffect"
$ws.Range("B3").Value = "Provide a short comparison of the two hooks"
$ws.Range("C3").Value = "Frontend"
$ws.Range("D3").Value = "text"
$ws.Range("H3").Value = "react,hooks"

# Row 4: Select all prime numbers -> Which of the following are fruits?
$ws.Range("A4").Value = "Which of the following are fruits?"
$ws.Range("B4").Value = "Select all that are fruits"
$ws.Range("C4").Value = "General Knowledge"
$ws.Range("D4").Value = "multi_choice"
$ws.Range("E4").Value = '["Apple","Carrot","Banana","Potato"]'
$ws.Range("F4").Value = '["Apple","Banana"]'
$ws.Range("H4").Value = "food,fruit"
